$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 10
$ws.Range("E5").Value = 10
$ws.Range("E6").Value = 10
$ws.Range("E7").Value = 3
$ws.Range("E12").Value = 3
$ws.Range("E13").Value = 3
$ws.Range("E14").Value = 3
$ws.Range("E15").Value = 4
$ws.Range("E16").Value = 4
$ws.Range("E18").Value = 4
$ws.Range("E19").Value = 4
$ws.Range("E20").Value = 4
$ws.Range("E22").Value = 4
$ws.Range("E23").Value = 4
$ws.Range("E24").Value = 4
$ws.Range("E25").Value = 4
$ws.Range("E26").Value = 4
$ws.Range("E27").Value = 12
$ws.Range("E28").Value = 12
$ws.Range("E30").Value = 12
$ws.Range("E32").Value = 12
$ws.Range("E33").Value = 12
$ws.Range("E34").Value = 12
$ws.Range("E36").Value = 12
$ws.Range("E38").Value = 12
$ws.Range("E41").Value = 2
$ws.Range("E42").Value = 2
$ws.Range("E43").Value = 2
$ws.Range("E50").Value = 2
$ws.Range("E51").Value = 2
$ws.Range("E53").Value = 2
$ws.Range("E54").Value = 2
$ws.Range("E55").Value = 8
$ws.Range("E56").Value = 8
$ws.Range("E57").Value = 8
$ws.Range("E60").Value = 8
$ws.Range("E62").Value = 8
$ws.Range("E66").Value = 8
$ws.Range("E69").Value = 13
$ws.Range("E70").Value = 13
$ws.Range("E72").Value = 13
$ws.Range("E73").Value = 13
$ws.Range("E74").Value = 13
$ws.Range("E75").Value = -1
$ws.Range("E76").Value = 10
$ws.Range("E79").Value = 10
$ws.Range("E80").Value = 10
$ws.Range("E81").Value = 3
$ws.Range("E82").Value = 10
$ws.Range("E84").Value = 1
$ws.Range("E85").Value = 1
$ws.Range("E86").Value = 1
$ws.Range("E87").Value = 1
$ws.Range("E88").Value = 1
$ws.Range("E89").Value = 1
$ws.Range("E91").Value = 1
$ws.Range("E92").Value = 1
$ws.Range("E97").Value = 9
$ws.Range("E98").Value = 9
$ws.Range("E100").Value = 9
$ws.Range("E101").Value = 9
$ws.Range("E102").Value = 9
$ws.Range("E103").Value = 9
$ws.Range("E104").Value = 9
$ws.Range("E105").Value = 9
$ws.Range("E111").Value = -1
$ws.Range("E114").Value = 11
$ws.Range("E115").Value = 11
$ws.Range("E116").Value = 8
$ws.Range("E118").Value = 11
$ws.Range("E119").Value = 0
$ws.Range("E121").Value = 0
$ws.Range("E122").Value = 0
$ws.Range("E123").Value = 0
$ws.Range("E126").Value = 0
$ws.Range("E128").Value = 0
$ws.Range("E129").Value = 0
$ws.Range("E130").Value = 7
$ws.Range("E131").Value = 7
$ws.Range("E132").Value = 7
$ws.Range("E133").Value = 7
$ws.Range("E134").Value = 7
$ws.Range("E135").Value = 7
$ws.Range("E136").Value = 7
$ws.Range("E137").Value = 7
$ws.Range("E138").Value = 7
$ws.Range("E139").Value = 7
$ws.Range("E140").Value = 7
$ws.Range("E141").Value = 7
$ws.Range("E142").Value = 5
$ws.Range("E147").Value = 5
$ws.Range("E149").Value = 5
$ws.Range("E151").Value = 5
$ws.Range("E153").Value = 13
$ws.Range("E154").Value = 13
$ws.Range("E155").Value = 6
$ws.Range("E156").Value = 6
$ws.Range("E158").Value = 13
$ws.Range("E159").Value = -1
$ws.Range("E160").Value = -1
$ws.Range("E161").Value = 13
$ws.Range("E162").Value = 3
$ws.Range("E163").Value = 3
$ws.Range("E171").Value = 3
$ws.Range("E172").Value = 10
$ws.Range("B176").Value = 1085.519200488409
$ws.Range("C176").Value = 1790.187777307487
$ws.Range("D176").Value = 105
$ws.Range("B177").Value = 2347.842460310387
$ws.Range("C177").Value = 1218.924660370858
$ws.Range("D177").Value = 84
$ws.Range("B178").Value = 80.25953621033399
$ws.Range("C178").Value = 539.3007741857029
$ws.Range("D178").Value = 68
$ws.Range("B179").Value = 878.5546239772101
$ws.Range("C179").Value = 2525.733703979399
$ws.Range("D179").Value = 114
$ws.Range("B180").Value = 2232.152399223904
$ws.Range("C180").Value = 2431.887445390206
$ws.Range("D180").Value = 72
$ws.Range("B181").Value = 1051.942866460077
$ws.Range("C181").Value = 528.8115188052595
$ws.Range("D181").Value = 104
$ws.Range("B182").Value = 445.1968076949909
$ws.Range("C182").Value = 2518.535861883758
$ws.Range("D182").Value = 114
$ws.Range("B183").Value = 2945.464993768458
$ws.Range("C183").Value = 1351.731251727125
$ws.Range("D183").Value = 105
$ws.Range("B184").Value = 381.5077754268682
$ws.Range("C184").Value = 982.9680528302903
$ws.Range("B185").Value = 1499.979726294803
$ws.Range("C185").Value = 1174.574709021907
$ws.Range("D185").Value = 75
$ws.Range("B186").Value = 1103.44024512642
$ws.Range("C186").Value = 2717.27020280983
$ws.Range("D186").Value = 150
$ws.Range("B187").Value = 636.1610031366665
$ws.Range("C187").Value = 782.3971833363092
$ws.Range("D187").Value = 67
$ws.Range("B188").Value = 1719.916974583251
$ws.Range("C188").Value = 2605.139153466509
$ws.Range("D188").Value = 41
$ws.Range("B189").Value = 426.4705458155561
$ws.Range("C189").Value = 2259.244686148031
$ws.Range("D189").Value = 83
